# Updated statistics: remove the outdated "ManuallyLoadCustodians" test suite entry
# (it has been replaced), and refresh the dependent summary formulas / shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the obsolete "ManuallyLoadCustodians" row entirely (A4:D4), including its
# formatting, so the row collapses back to an untouched state.
$ws.Range("A4:D4").Clear() | Out-Null

# The "Not yet re-implemented" status note in E4 no longer applies; clear its value
# but keep the cell's existing border/style.
$ws.Range("E4").ClearContents() | Out-Null

# Reset the view so the sheet opens scrolled to the top-left (no frozen/scrolled
# "topLeftCell"), with the now-empty former data row selected.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("A4:E4").Select() | Out-Null
